# The commit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: the deck's main theme ("Integral" / Red Violet,
# wired to the slide master + presentation.xml) becomes the stock
# "Office Theme" (blue accent1 etc.), while the theme that was only
# referenced by the notes master becomes "Integral" / Red Violet.
#
# Both theme parts already share an identical <a:fontScheme> and
# <a:fmtScheme> - the only real difference between them is the 12
# <a:clrScheme> colours (and the cosmetic name="" attributes), so the
# edit is reproduced here by re-pointing the live theme's colour
# scheme at the "Office" palette.

function Get-RgbLong([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# Target palette == the current "Office Theme" colours that live in
# ppt/theme/theme1.xml (clrScheme name="Office") before the swap.
$officeColors = @(
    @(0x00, 0x00, 0x00), # 1  dk1
    @(0xFF, 0xFF, 0xFF), # 2  lt1
    @(0x44, 0x54, 0x6A), # 3  dk2
    @(0xE7, 0xE6, 0xE6), # 4  lt2
    @(0x5B, 0x9B, 0xD5), # 5  accent1
    @(0xED, 0x7D, 0x31), # 6  accent2
    @(0xA5, 0xA5, 0xA5), # 7  accent3
    @(0xFF, 0xC0, 0x00), # 8  accent4
    @(0x44, 0x72, 0xC4), # 9  accent5
    @(0x70, 0xAD, 0x47), # 10 accent6
    @(0x05, 0x63, 0xC1), # 11 hlink
    @(0x95, 0x4F, 0x72)  # 12 folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $rgb = $officeColors[$i]
    $themeColors.Colors($i + 1).RGB = Get-RgbLong $rgb[0] $rgb[1] $rgb[2]
}

# Best-effort cosmetic rename (some hosts keep these read-only/stubbed).
try { $master.Theme.Name = "Office Theme" } catch {}
try { $themeColors.Name = "Office" } catch {}
